$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled sheet name "Oppertunity" -> "Opportunity"
$ws.Name = "Opportunity"

# Move/update the active selection from I1 to C7
$ws.Range("C7").Select()
